$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$words = @(
    "reizen",
    "knarren",
    "kaufen",
    "schenken",
    "wirken",
    "saufen",
    "werfen",
    "fischen",
    "wohnen",
    "bellen",
    "rufen",
    "boxen",
    "zünden",
    "geben",
    "bauen",
    "leihen",
    "wundern",
    "schmecken",
    "schalten",
    "kümmern",
    "kriegen",
    "freuen",
    "morden",
    "kosten",
    "schulden",
    "schwingen",
    "wagen",
    "rasen",
    "runden",
    "brauen",
    "fühlen",
    "münzen",
    "gnaden",
    "gelten",
    "decken",
    "fehlen",
    "dienen",
    "bergen",
    "folgen",
    "achten",
    "stehlen",
    "spannen",
    "fällen",
    "feiern",
    "lügen",
    "starten",
    "passen",
    "warnen",
    "kichern",
    "gründen",
    "loben",
    "siegen",
    "leiden",
    "brauchen",
    "knien",
    "wehen",
    "fallen",
    "malen",
    "küssen",
    "bitten",
    "scheitern",
    "stören",
    "kehren",
    "wehtun",
    "leugnen",
    "sterben",
    "schämen",
    "mögen",
    "spinnen",
    "fragen",
    "irren",
    "biegen",
    "sinken",
    "quälen",
    "spielen",
    "arten",
    "schütteln",
    "hauen",
    "pflanzen",
    "jubeln",
    "wüten",
    "führen",
    "süßen",
    "wenden",
    "flehen",
    "kratzen",
    "posten",
    "eignen",
    "grüßen",
    "zahlen",
    "liegen",
    "scheinen",
    "drücken",
    "dauern",
    "träumen",
    "knurren",
    "sperren",
    "schwächen",
    "spüren",
    "füttern",
    "trauen",
    "locken",
    "streichen",
    "holen",
    "atmen",
    "schlucken",
    "treiben",
    "regeln",
    "heulen",
    "stecken",
    "schwören",
    "stärken",
    "planen",
    "flüchten",
    "betteln",
    "graben",
    "backen",
    "zählen",
    "fahren",
    "sichern",
    "ärgern",
    "lieben",
    "stellen",
    "heben",
    "fordern",
    "äußern",
    "lesen",
    "tollen",
    "erben",
    "schrecken",
    "bluten",
    "ruhen",
    "greifen",
    "machen",
    "weichen",
    "heißen",
    "hören",
    "helfen",
    "spenden",
    "töten",
    "ändern",
    "haben",
    "stechen",
    "klingen",
    "drehen",
    "zeichnen",
    "zögern",
    "wachsen",
    "starren",
    "schreiten",
    "pfeifen",
    "jagen",
    "reisen",
    "zeugen",
    "melden",
    "sprengen",
    "bleiben",
    "formen",
    "altern",
    "fangen",
    "filmen",
    "stammen",
    "klingeln",
    "flüstern",
    "herrschen",
    "lockern",
    "klettern",
    "dringen",
    "enden",
    "werden",
    "sammeln",
    "heilen",
    "wollen",
    "räumen",
    "trennen",
    "trösten",
    "zielen",
    "zeigen",
    "danken",
    "liefern",
    "binden",
    "suchen",
    "tropfen",
    "platzen",
    "ehren",
    "baden",
    "sorgen",
    "sagen",
    "fügen",
    "mauern",
    "seufzen",
    "fließen"
)
for ($i = 0; $i -lt $words.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $words[$i]
}